{"js": "// no-op test\ncontext.document.body.load(\"text\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$d.Content.Text | Out-Null\n"}
